$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2024-07-03 Wednesday" "2024-07-04 Thursday"

Replace-Text "684÷5=" "319÷6="
Replace-Text "423÷6=" "164÷3="
Replace-Text "867÷7=" "315÷2="
Replace-Text "402÷4=" "436÷7="
Replace-Text "707÷5=" "545÷2="
Replace-Text "200÷4=" "457÷9="
Replace-Text "754÷6=" "185÷6="
Replace-Text "653÷7=" "605÷5="
Replace-Text "726÷9=" "687÷8="
Replace-Text "985÷5=" "167÷7="
Replace-Text "721÷3=" "779÷5="
Replace-Text "978÷3=" "694÷6="
Replace-Text "549÷4=" "544÷7="
Replace-Text "218÷9=" "765÷6="
Replace-Text "755÷8=" "556÷7="
Replace-Text "856÷7=" "302÷6="
Replace-Text "372÷3=" "265÷9="
Replace-Text "989÷5=" "410÷2="
Replace-Text "899÷8=" "506÷5="
Replace-Text "915÷7=" "810÷2="
Replace-Text "389÷6=" "134÷6="
Replace-Text "376÷5=" "178÷9="
Replace-Text "860÷3=" "397÷7="
Replace-Text "896÷9=" "982÷2="
Replace-Text "594÷9=" "967÷8="
